$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Ecoli short reads": bump a couple of read-count values and
# move the view/selection.
# ---------------------------------------------------------------------
$wsEcoliShort = $wb.Worksheets.Item("Ecoli short reads")
$wsEcoliShort.Range("B3").Value = 332007
$wsEcoliShort.Range("D34").Select()

# ---------------------------------------------------------------------
# Sheet "Sevim Real Data": add a new header row above the existing
# table (pushing everything else down one row) with two new labels,
# narrow column B, and move the selection.
# ---------------------------------------------------------------------
$wsSevim = $wb.Worksheets.Item("Sevim Real Data")
$wsSevim.Rows("1:1").Insert()
$wsSevim.Range("B1").Value = "Percent of Bases Assigned to each Reference"
$wsSevim.Range("E1").Value = "Basses Assigned to each Reference"
$wsSevim.Columns("B").ColumnWidth = 15.85
$wsSevim.Range("E10").Select()

# ---------------------------------------------------------------------
# Sheet "Time and Memory": add two subtotal formulas and move the
# selection.
# ---------------------------------------------------------------------
$wsTime = $wb.Worksheets.Item("Time and Memory")
$wsTime.Range("C11").Formula = "=C6+C5"
$wsTime.Range("C22").Formula = "=C17+C16"
$wsTime.Range("H6").Select()

# ---------------------------------------------------------------------
# Sheet "Simulated Data": correct the H34/I34 read counts and move the
# selection. Re-activate this sheet last so it stays the active tab.
# ---------------------------------------------------------------------
$wsSim = $wb.Worksheets.Item("Simulated Data")
$wsSim.Range("H34").Value = 642373
$wsSim.Range("I34").Value = 0
$wsSim.Activate()
$wsSim.Range("I37").Select()
